# Updates cryptos list figures (price + 1h volume change) to match the
# upstream GitHub Actions refresh commit. Price cells that look purely
# numeric (single decimal point, e.g. "1.060") are written with a leading
# apostrophe so Excel keeps them as text (preserving trailing zeros etc.)
# exactly like the source workbook, instead of silently coercing them to
# numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "20.529.34"
$ws.Range("E2").Value = "  +2.55%  "
$ws.Range("D3").Value = "1.471.23"
$ws.Range("E3").Value = "  +3.55%  "
$ws.Range("D4").Value = "'1.008"
$ws.Range("E4").Value = "  +0.74%  "
$ws.Range("D5").Value = "'0.9442"
$ws.Range("E5").Value = "  -5.61%  "
$ws.Range("D6").Value = "'281.09"
$ws.Range("E6").Value = "  +2.46%  "
$ws.Range("D7").Value = "'0.3716"
$ws.Range("E7").Value = "  -0.19%  "
$ws.Range("E8").Value = "  +3.29%  "
$ws.Range("D9").Value = "'41.44"
$ws.Range("E9").Value = "  +3.81%  "
$ws.Range("D10").Value = "'1.060"
$ws.Range("E10").Value = "  +4.72%  "
$ws.Range("D11").Value = "'0.06683"
$ws.Range("E11").Value = "  +1.23%  "
$ws.Range("D12").Value = "'1.003"
$ws.Range("E12").Value = "  +0.31%  "
$ws.Range("D13").Value = "'5.602"
$ws.Range("E13").Value = "  +3.38%  "
$ws.Range("D14").Value = "'18.31"
$ws.Range("E14").Value = "  +6.71%  "
$ws.Range("D15").Value = "'6.239"
$ws.Range("E15").Value = "  +0.87%  "
$ws.Range("D16").Value = "1.476.24"
$ws.Range("E16").Value = "  +3.79%  "
$ws.Range("E17").Value = "  +2.49%  "
$ws.Range("D18").Value = "'0.9368"
$ws.Range("E18").Value = "  -6.34%  "
$ws.Range("D19").Value = "'0.05748"
$ws.Range("E19").Value = "  -0.92%  "
$ws.Range("D20").Value = "'72.38"
$ws.Range("E20").Value = "  -3.04%  "
$ws.Range("D21").Value = "'5.702"
$ws.Range("E21").Value = "  +1.10%  "
$ws.Range("E22").Value = "  +1.91%  "
$ws.Range("E23").Value = "  +2.17%  "
$ws.Range("D24").Value = "'2.271"
$ws.Range("E24").Value = "  -2.47%  "
$ws.Range("D25").Value = "20.798.99"
$ws.Range("E25").Value = "  +3.87%  "
$ws.Range("D26").Value = "'2.310"
$ws.Range("E26").Value = "  +0.58%  "
$ws.Range("D27").Value = "'138.11"
$ws.Range("E27").Value = "  -0.79%  "
$ws.Range("E28").Value = "  +4.15%  "
$ws.Range("D29").Value = "1.637.16"
$ws.Range("E29").Value = "  +3.59%  "
$ws.Range("D30").Value = "'113.90"
$ws.Range("E30").Value = "  +4.21%  "
$ws.Range("D31").Value = "'3.952"
$ws.Range("E31").Value = "  +3.71%  "
$ws.Range("D32").Value = "'5.321"
$ws.Range("E32").Value = "  -1.76%  "
$ws.Range("D33").Value = "'0.8514"
$ws.Range("E33").Value = "  -4.44%  "
$ws.Range("D34").Value = "'1.602"
$ws.Range("E34").Value = "  +25.59%  "
$ws.Range("E35").Value = "  +1.22%  "
$ws.Range("D36").Value = "'0.06062"
$ws.Range("E36").Value = "  +5.45%  "
$ws.Range("D37").Value = "'4.943"
$ws.Range("E37").Value = "  +3.00%  "
$ws.Range("D38").Value = "'10.74"
$ws.Range("E38").Value = "  -4.91%  "
$ws.Range("D39").Value = "'0.02074"
$ws.Range("E39").Value = "  +1.64%  "
$ws.Range("D40").Value = "'1.123"
$ws.Range("E40").Value = "  +3.76%  "
$ws.Range("B41").Value = "Algorand"
$ws.Range("C41").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D41").Value = "'0.1908"
$ws.Range("E41").Value = "  -1.09%  "
$ws.Range("B42").Value = "Frax"
$ws.Range("C42").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D42").Value = "'0.9571"
$ws.Range("E42").Value = "  -4.27%  "
$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D43").Value = "'7.529"
$ws.Range("E43").Value = "  -10.77%  "
$ws.Range("D44").Value = "'0.5418"
$ws.Range("E44").Value = "  +1.62%  "
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").Value = "'12.55"
$ws.Range("E45").Value = "  +2.25%  "
$ws.Range("B46").Value = "PancakeSwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D46").Value = "'3.584"
$ws.Range("E46").Value = "  +1.25%  "
$ws.Range("D47").Value = "'122.31"
$ws.Range("E47").Value = "  +11.36%  "
$ws.Range("E48").Value = "  +3.85%  "
$ws.Range("E49").Value = "  +1.46%  "
$ws.Range("E50").Value = "  +4.34%  "
$ws.Range("E51").Value = "  -0.42%  "
